# Insercion y exportacion de productos
# Update the "services" import template:
#  - header row (row 7) is replaced with the new product-import column names
#  - the two non-red header cells (B7, F7) get an explicit black font color

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for row 7 (same bold styling already present per cell)
$ws.Range("A7").Value = "Nombre"
$ws.Range("B7").Value = "Cantidad"
$ws.Range("C7").Value = "Codigo"
$ws.Range("D7").Value = "id_categoria"
$ws.Range("E7").Value = "fotografia"
$ws.Range("F7").Value = "descripcion"
$ws.Range("G7").Value = "pv"
$ws.Range("H7").Value = "pc"
$ws.Range("I7").Value = "id_proveedor"
$ws.Range("J7").Value = "estado"

# B7 / F7 keep bold, but now carry an explicit black font color
$ws.Range("B7").Font.Color = 0
$ws.Range("F7").Font.Color = 0

# Cursor ends up parked on L9, like in the saved workbook
$ws.Range("L9").Select()
